# Reverts the "Hotfix for shape setting, fixed column in excel sheet loaded" change:
#  - Custom Shape? value (P4 on Sheet1) goes back from "False" to "True"
#  - The sheet's active selection moves from P4 to J13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restore P4 ("Custom Shape?" row) to "True" ---
# Simply assigning the string "True"/"False" gets auto-coerced to a boolean
# cell by Excel. To faithfully restore the original shared-string cell
# (t="s", same style, reusing the existing "True" shared string instead of
# creating a boolean or a brand new string entry) we copy an existing cell
# that already holds the "True" text with the same style (G3) onto P4.
$trueCell = $ws.Cells.Item(3, 7)   # G3 => "True", style s="1"
$targetCell = $ws.Cells.Item(4, 16) # P4
$trueCell.Copy($targetCell)

# --- Update the selected/active cell in the sheet view ---
$ws.Range("J13").Select()
